$d = $word.ActiveDocument

$LDQ = [char]8220   # “
$RDQ = [char]8221   # ”
$RSQ = [char]8217   # '  (right single quote / apostrophe)

# ---------------------------------------------------------------------
# 1. Add <w:noProof/> to the run holding the ER Diagram picture.
# ---------------------------------------------------------------------
$shp = $d.InlineShapes.Item(1)
$shp.Range.NoProofing = 1

# ---------------------------------------------------------------------
# Helper: locate the paragraph whose text starts with a given prefix.
# ---------------------------------------------------------------------
function Get-ParagraphByPrefix($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "$prefix*") {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 2. Relation Schema: Calendar ( ... )
#    ", owner_sin, date_start, date_end, price)"
#    -> ", <u>date_start</u>, <u>date_end</u>, price)"
# ---------------------------------------------------------------------
$p = Get-ParagraphByPrefix("Calendar (")
$r = $p.Range
$r.Find.Execute("owner_sin, ") | Out-Null
$r.Text = ""

$r2 = $p.Range
$r2.Find.Execute("date_start") | Out-Null
$r2.Font.Underline = 1

$r3 = $p.Range
$r3.Find.Execute("date_end") | Out-Null
$r3.Font.Underline = 1

# ---------------------------------------------------------------------
# 3. Relation Schema: Listing_Comments ( ... )
#    ", comments" -> ", <u>timestamp</u>, comments"
# ---------------------------------------------------------------------
$p = Get-ParagraphByPrefix("Listing_Comments (")
$r = $p.Range
$r.Find.Execute(", comments") | Out-Null
$r.Text = ", timestamp, comments"

$r2 = $p.Range
$r2.Find.Execute("timestamp") | Out-Null
$r2.Font.Underline = 1

# ---------------------------------------------------------------------
# 4. Relation Schema: Rent_History ( ... )
#    ", sin, type" -> ", <u>sin</u>, type"
# ---------------------------------------------------------------------
$p = Get-ParagraphByPrefix("Rent_History (")
$r = $p.Range
$r.Find.Execute("sin") | Out-Null
$r.Font.Underline = 1

# ---------------------------------------------------------------------
# 5. Relation Schema: Ownership ( ... )
#    ", listno" -> ", <u>listno</u>"
# ---------------------------------------------------------------------
$p = Get-ParagraphByPrefix("Ownership (")
$r = $p.Range
$r.Find.Execute("listno") | Out-Null
$r.Font.Underline = 1

# ---------------------------------------------------------------------
# 6. "The text interface will appear on the console to the right. ..."
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("console to the right", $true, $false, $false, $false, $false, $true, 1, $false, "console in the IDE", 2) | Out-Null

$r = $d.Content
$needle = "a number or " + $LDQ + "y" + $RDQ + " and hitting"
$r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "a number and hitting", 2) | Out-Null

# ---------------------------------------------------------------------
# 7. "System limitations are about installed software. ..."
#    "Java's Connector-J"  -> "SQL's Connector/J"
#    "Java/MySQL"          -> "Eclipse/Java/SQL"
# ---------------------------------------------------------------------
$r = $d.Content
$needle = "Java" + $RSQ + "s Connector-J"
$repl = "SQL" + $RSQ + "s Connector/J"
$r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $repl, 2) | Out-Null

$r = $d.Content
$r.Find.Execute("cannot run Java/MySQL", $true, $false, $false, $false, $false, $true, 1, $false, "cannot run Eclipse/Java/SQL", 2) | Out-Null

# ---------------------------------------------------------------------
# 8. Final paragraph restructuring:
#    - split after "... catch statements."
#    - add a new sentence about the local website idea
#    - start a new paragraph "Regarding SQL, I feel alright..."
#    - append a closing sentence about SQL query accuracy
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(" In regards to the SQL side of things, ") | Out-Null
$r.Text = "`r"

$r = $d.Content
$r.Find.Execute("catch statements.") | Out-Null
$r.Collapse(0)
$r.InsertAfter(" This lack of focus on Java also made the interface difficult to progress quickly. A local website attached to a database is one area this could be improved upon. ")

$r = $d.Content
$r.Find.Execute("I feel alright") | Out-Null
$r.Collapse(1)
$r.InsertBefore("Regarding SQL, ")

$r = $d.Content
$r.Find.Execute("direct SQL answer.") | Out-Null
$r.Collapse(0)
$r.InsertAfter(" The SQL queries seem to work for simple questions.")
